$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = Extension -> Constraint(s) column AJ
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 = Extension.id -> Type(s) column K
$elements.Range("K3").Value = "string
"

# Row 6 = Extension.value[x] -> Definition column M
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
